# Act 2-2 dialog stuff
# - Update the "nefarious goblins" line (row 85, B column) to new wording.
# - Fill in the previously-blank placeholder rows 90-94 with new
#   "newton_second_law_2_dlg_*" key/value pairs (cannon introduction scene).
# - Move the active selection / scroll position to reflect where the
#   author ended up editing (B92).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 85: reword the existing goblins dialogue line ---
$ws.Range("B85").Value = "The nefarious goblins have appeared out of thin air! They are surely up to no good. Push them off the cliff using the wheel."
$ws.Range("B85").VerticalAlignment = -4108   # xlCenter -> matches style s="2"

# --- Row 90: newton_second_law_2_dlg_1 ---
$ws.Range("A90").Value = "newton_second_law_2_dlg_1"
$ws.Range("B90").Value = "Hark! More goblins have appeared! This time, they have positioned themselves at different heights."
$ws.Range("B90").Style = "Normal"

# --- Row 91: newton_second_law_2_dlg_2 ---
$ws.Range("A91").Value = "newton_second_law_2_dlg_2"
$ws.Range("B91").Value = "But fear not, we have the very tool to get the job done."
$ws.Range("B91").VerticalAlignment = -4108   # xlCenter -> matches style s="2"

# --- Row 92: newton_second_law_2_dlg_3 ---
$ws.Range("A92").Value = "newton_second_law_2_dlg_3"
$ws.Range("B92").Value = "Let us bring forth the mighty cannon to vanquish these vermin!"
$ws.Range("B92").VerticalAlignment = -4108   # xlCenter -> matches style s="2"

# --- Row 93: newton_second_law_2_dlg_4 ---
$ws.Range("A93").Value = "newton_second_law_2_dlg_4"
$ws.Range("B93").Value = "In this scenario, we are applying force to a cannonball with explosion."
$ws.Range("B93").Style = "Normal"

# --- Row 94: newton_second_law_2_dlg_5 ---
$ws.Range("A94").Value = "newton_second_law_2_dlg_5"
$ws.Range("B94").Value = "This short burst of force will allow the cannonball to accelerate within a fraction of a second to reach high velocity."
$ws.Range("B94").Style = "Normal"

# --- View state: scroll/selection moved down to the newly-edited area ---
$win = $excel.ActiveWindow
$win.ScrollRow = 73
$win.ScrollColumn = 1
$ws.Range("B92").Select()
